{"js": "// Replace \"e-billing\" with \"billing\" everywhere in the document body.\n// (The source doc used \"e-billing\" consistently; the edit drops the \"e-\" prefix,\n// e.g. \"e-billing account number\" -> \"billing account number\".)\nconst body = context.document.body;\nconst results = body.search(\"e-billing\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"billing\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace \"e-billing\" with \"billing\" throughout the document.\n# (The document consistently used \"e-billing\"; the edit drops the \"e-\" prefix,\n# e.g. \"e-billing account number\" -> \"billing account number\".)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"e-billing\"\n$find.Replacement.Text = \"billing\"\n$find.Forward = $true\n$find.Wrap = 1        # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $false,\n    $false,\n    $find.Forward,\n    $find.Wrap,\n    $find.Format,\n    $find.Replacement.Text,\n    2                  # wdReplaceAll\n)\n"}
